$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$th = $sm.Theme
$cs = $th.ThemeColorScheme
try {
    $cs.Load("Office")
    Write-Output "load ok"
} catch {
    Write-Output "ERR: $_"
}
